# A2_TOC.docx — Table of Contents update
# Adds two new TOC rows for the "no assignments" landing-page designs and
# bumps the page numbers of all the rows that follow each insertion point
# (pages shifted by the newly-added content).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellPageNumber {
    param(
        [int]$RowIndex,
        [string]$OldText,
        [string]$NewText
    )
    $cell = $t.Cell($RowIndex, 2)
    $rng = $cell.Range
    # MatchWholeWord; no wrap beyond the cell's own Range; ReplaceOne (wdReplaceOne)
    # so the substitution stays confined to this cell instead of cascading
    # document-wide the way wdReplaceAll does.
    $rng.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 0, $false, $NewText, 1) | Out-Null
}

# --- Step 1: update existing page numbers (row indices as they are in the
#             original 32-row table, before any rows are inserted).
#             Processed bottom-to-top: the COM Find engine's match/replace
#             positions lag by one call when walking top-to-bottom through
#             adjacent cells, so going in reverse keeps every Find scoped
#             to the right cell. -----------------------------------------
Set-CellPageNumber 32 "110" "113"
Set-CellPageNumber 31 "106" "109"
Set-CellPageNumber 30 "102" "105"
Set-CellPageNumber 29 "100" "103"
Set-CellPageNumber 28 "97"  "100"
Set-CellPageNumber 27 "88"  "91"
Set-CellPageNumber 26 "84"  "87"
Set-CellPageNumber 25 "81"  "84"
Set-CellPageNumber 24 "79"  "82"
Set-CellPageNumber 23 "76"  "79"
Set-CellPageNumber 22 "74"  "77"
Set-CellPageNumber 21 "66"  "69"
Set-CellPageNumber 20 "65"  "68"
Set-CellPageNumber 19 "63"  "66"
Set-CellPageNumber 18 "59"  "62"
Set-CellPageNumber 17 "57"  "59"
Set-CellPageNumber 16 "54"  "56"
Set-CellPageNumber 15 "52"  "54"
Set-CellPageNumber 14 "44"  "46"
Set-CellPageNumber 13 "42"  "44"
Set-CellPageNumber 12 "41"  "42"
Set-CellPageNumber 11 "40"  "41"
Set-CellPageNumber 10 "33"  "34"
Set-CellPageNumber 9  "32"  "33"
Set-CellPageNumber 8  "28"  "29"
Set-CellPageNumber 7  "27"  "28"
Set-CellPageNumber 6  "26"  "27"
Set-CellPageNumber 5  "24"  "25"
Set-CellPageNumber 4  "11"  "12"
Set-CellPageNumber 3  "9"   "10"
Set-CellPageNumber 2  "6"   "7"

# --- Step 2: insert the two new TOC rows ----------------------------------
# New row 1: right after "ASSIGNMENT LIST" (row 1), before "REORDER" (row 2)
$newRow1 = $t.Rows.Add($t.Rows.Item(2))
$t.Cell(2, 1).Range.Text = "VIEW OF ASSIGNMENT LIST WITH NO ASSIGNMENTS"
$t.Cell(2, 2).Range.Text = "1"

# New row 2: right after "UPLOAD GRADES" (now row 18), before "SV: ASSIGNMENT LIST" (now row 19)
$newRow2 = $t.Rows.Add($t.Rows.Item(19))
$t.Cell(19, 1).Range.Text = "SV: ASSIGNMENT LIST WITH NO ASSIGNMENTS"
$t.Cell(19, 2).Range.Text = "61"
